$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update every row in column C (the "Förändrad" date column) from 45204 to 45205.
# Only rows that currently hold the old value are touched, keeping the edit
# scoped exactly to what changed in the source workbook.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45205
    }
}
